# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E) figures
# to match the latest pull from the data feed (GitHub Actions run of
# Mon Dec 18 13:55:13 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.263.01"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.153.19"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("D15").Value = "2.472.23"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "2.134.74"
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.781"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.92%  "
$ws.Range("D19").Value = "41.194.45"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("E20").Value = "  -4.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("E22").Value = "  -6.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "225.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("E27").Value = "  -7.16%  "
$ws.Range("E28").Value = "  -9.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.46%  "
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.09%  "
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.39%  "
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.77%  "
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("E51").Value = "  -3.07%  "
